$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-491 all contain serial date 45206 (2023-10-07),
# which should be updated to 45208 (2023-10-09).
$ws.Range("C2:C491").Value = 45208
